# Generate Report for Handoff
# Adds two new "handed off" rows (for the files 6b9ff258-...md and
# 965104be-...md) to the Overview sheet and to each locale sheet
# (zh-cn / de-de), mirroring the existing "Ready for handoff" rows.

$wb = $excel.ActiveWorkbook

# ---- shared constants -------------------------------------------------
$file1 = "6b9ff258-5cf7-40c6-93dc-10784c1d9a46"
$file2 = "965104be-0ef5-4edb-82b1-facb0f37c968"

$hash1 = "0abd8208a0eb44c32d9f52395849381cb7361d5f"
$hash2 = "a05aeaf165ae5501d28d8aece37eedefb60075ad"

$overviewDate = "2016-25-11 16:25:17"

$zhDate1 = "2016-03-11 16:25:14"
$zhDate2 = "2016-03-11 16:25:14"

$deDate1 = "2016-03-11 16:25:17"
$deDate2 = "2016-03-11 16:25:17"

$noHandback = "0001-01-01 00:00:00"

$linkColor = 15570276   # RGB(0x64, 0x95, 0xED) -> matches existing HyperLink font color
$dateFormat = "yyyy-mm-dd HH:mm:ss"

function Style-AsLink($range) {
    $range.Font.Name = "Calibri"
    $range.Font.Size = 11
    $range.Font.Underline = 2
    $range.Font.Color = $linkColor
}

# ========================================================================
# Sheet "Overview"
# ========================================================================
$ws = $wb.Worksheets.Item("Overview")

# --- row 4 : 6b9ff258-... ---
$ws.Range("A4").Value = "$file1.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = $overviewDate

Style-AsLink $ws.Range("A4")
$ws.Hyperlinks.Add(
    $ws.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/main/e2e/$file1.md",
    [Type]::Missing,
    [Type]::Missing,
    "$file1.md") | Out-Null

# --- row 5 : 965104be-... ---
$ws.Range("A5").Value = "$file2.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = $overviewDate

Style-AsLink $ws.Range("A5")
$ws.Hyperlinks.Add(
    $ws.Range("A5"),
    "https://github.com/OpenLocalizationTest/oltest/blob/main/e2e/$file2.md",
    [Type]::Missing,
    [Type]::Missing,
    "$file2.md") | Out-Null

# ========================================================================
# Locale sheets: zh-cn / de-de
# ========================================================================
$locales = @(
    @{ Name = "zh-cn"; Date1 = $zhDate1; Date2 = $zhDate2 },
    @{ Name = "de-de"; Date1 = $deDate1; Date2 = $deDate2 }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Name)
    $lc = $loc.Name

    $xlf1 = "$file1.$hash1.$lc.xlf"
    $xlf2 = "$file2.$hash2.$lc.xlf"

    # --- row 4 : 6b9ff258-... ---
    $ws.Range("A4").Value = "$file1.md"
    $ws.Range("B4").Value = ".md"
    $ws.Range("C4").Value = "Ready for handoff"
    $ws.Range("D4").Value = $xlf1
    $ws.Range("E4").Value = $loc.Date1
    $ws.Range("E4").NumberFormat = $dateFormat
    $ws.Range("H4").Value = $noHandback
    $ws.Range("I4").Value = "Include"

    Style-AsLink $ws.Range("A4")
    Style-AsLink $ws.Range("B4")
    Style-AsLink $ws.Range("D4")

    $ws.Hyperlinks.Add(
        $ws.Range("A4"),
        "https://github.com/OpenLocalizationTest/oltest/blob/main/e2e/$file1.md",
        [Type]::Missing,
        [Type]::Missing,
        "$file1.md") | Out-Null
    $ws.Hyperlinks.Add(
        $ws.Range("B4"),
        "https://github.com/OpenLocalizationTest/oltest/blob/main/e2e/$file1.md",
        [Type]::Missing,
        [Type]::Missing,
        ".md") | Out-Null
    $ws.Hyperlinks.Add(
        $ws.Range("D4"),
        "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/main/ol-handoff/OpenLocalizationTestOrg/oltest.$lc/ci/ht/$xlf1",
        [Type]::Missing,
        [Type]::Missing,
        $xlf1) | Out-Null

    # --- row 5 : 965104be-... ---
    $ws.Range("A5").Value = "$file2.md"
    $ws.Range("B5").Value = ".md"
    $ws.Range("C5").Value = "Ready for handoff"
    $ws.Range("D5").Value = $xlf2
    $ws.Range("E5").Value = $loc.Date2
    $ws.Range("E5").NumberFormat = $dateFormat
    $ws.Range("H5").Value = $noHandback
    $ws.Range("I5").Value = "Include"

    Style-AsLink $ws.Range("A5")
    Style-AsLink $ws.Range("B5")
    Style-AsLink $ws.Range("D5")

    $ws.Hyperlinks.Add(
        $ws.Range("A5"),
        "https://github.com/OpenLocalizationTest/oltest/blob/main/e2e/$file2.md",
        [Type]::Missing,
        [Type]::Missing,
        "$file2.md") | Out-Null
    $ws.Hyperlinks.Add(
        $ws.Range("B5"),
        "https://github.com/OpenLocalizationTest/oltest/blob/main/e2e/$file2.md",
        [Type]::Missing,
        [Type]::Missing,
        ".md") | Out-Null
    $ws.Hyperlinks.Add(
        $ws.Range("D5"),
        "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/main/ol-handoff/OpenLocalizationTestOrg/oltest.$lc/ci/ht/$xlf2",
        [Type]::Missing,
        [Type]::Missing,
        $xlf2) | Out-Null
}
